# magmaCooler scene: add two new localization rows ("rock_result" /
# "continue") right after the grain-size rows (old row 28), pushing the
# existing rock-type rows down by two. Mirrors the OOXML diff where rows
# 29-98 become rows 31-100 and two fresh rows are inserted at 29-30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 29-30; everything from the old row 29 onward
# (including the trailing numeric C-column weights) shifts down to 31-100.
$ws.Rows("29:30").Insert()

# Populate the new rows. Write row 30 ("continue") before row 29
# ("rock_result") so the new shared-string entries land in the same order
# as the target workbook: continue, CONTINUE, rock_result, ROCK RESULT.
$ws.Range("A30").Value = "continue"
$ws.Range("B30").Value = "CONTINUE"

$ws.Range("A29").Value = "rock_result"
$ws.Range("B29").Value = "ROCK RESULT"

# Match the author's updated viewport/selection (scrolled down a bit,
# new active cell on the freshly inserted row).
$win = $excel.ActiveWindow
$win.ScrollRow = 18
$win.ScrollColumn = 1
$ws.Range("A29").Select()
